# Insert a new weekly price record at row 118 for "Terminal Hortofrutícola
# Agro Chillán - Mango" (2023-09-07, Brasil origin), pushing the existing
# records (previously rows 118-189) down to rows 119-190.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 118; everything below shifts down one row.
$ws.Rows.Item(118).Insert()

# Seed the new row with a copy of the (now shifted) row directly below it,
# so the shared/static columns (Mercado, Región, Producto, Calidad, Unidad,
# etc.) carry over, matching how the rest of this table's rows look.
$ws.Range("A119:T119").Copy($ws.Range("A118:T118"))

# Overwrite the fields that differ for this new record.
$ws.Range("D118").Value = 45176
$ws.Range("M118").Value = 40
$ws.Range("N118").Value = 12000
$ws.Range("O118").Value = 12000
$ws.Range("P118").Value = 12000
$ws.Range("R118").Value = "Brasil"
$ws.Range("S118").Value = 3000
